# Apply crypto price/volume updates generated on Mon Jun 17 15:36:20 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($Cell, $Text) {
    $rng = $ws.Range($Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "65.423.32"
Set-TextValue "E2" "  -1.78%  "
Set-TextValue "D3" "3.509.03"
Set-TextValue "E3" "  -2.17%  "
Set-TextValue "E4" "  -0.02%  "
Set-TextValue "D5" "599.55"
Set-TextValue "E5" "  -1.57%  "
Set-TextValue "D6" "142.66"
Set-TextValue "E6" "  -3.12%  "
Set-TextValue "D7" "3.511.86"
Set-TextValue "E7" "  -2.05%  "
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.20%  "
Set-TextValue "E9" "  +5.29%  "
Set-TextValue "E10" "  -2.93%  "
Set-TextValue "E11" "  -2.22%  "
Set-TextValue "D12" "0.403"
Set-TextValue "D13" "4.112.50"
Set-TextValue "E13" "  -2.00%  "
Set-TextValue "D14" "0.0000195"
Set-TextValue "E14" "  -6.50%  "
Set-TextValue "D15" "28.25"
Set-TextValue "E15" "  -6.00%  "
Set-TextValue "D16" "3.506.68"
Set-TextValue "E16" "  -2.22%  "
Set-TextValue "E17" "  +1.40%  "
Set-TextValue "D18" "65.381.79"
Set-TextValue "E18" "  -1.95%  "
Set-TextValue "D19" "10.77"
Set-TextValue "E19" "  -5.64%  "
Set-TextValue "E20" "  -2.45%  "
Set-TextValue "D21" "14.34"
Set-TextValue "E21" "  -5.02%  "
Set-TextValue "D22" "418.38"
Set-TextValue "E22" "  -3.21%  "
Set-TextValue "D23" "0.593"
Set-TextValue "E23" "  -5.05%  "
Set-TextValue "D24" "76.93"
Set-TextValue "E24" "  -2.66%  "
Set-TextValue "D25" "3.655.25"
Set-TextValue "E25" "  -2.02%  "
Set-TextValue "E26" "  -0.02%  "
Set-TextValue "E27" "  -5.85%  "
Set-TextValue "E28" "  -3.01%  "
Set-TextValue "B29" "InternetComputer(DFINITY)"
Set-TextValue "C29" "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue "D29" "8.87"
Set-TextValue "E29" "  -4.72%  "
Set-TextValue "B30" "RenderToken"
Set-TextValue "C30" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D30" "7.71"
Set-TextValue "E30" "  -5.84%  "
Set-TextValue "D31" "0.999"
Set-TextValue "E31" "  -0.10%  "
Set-TextValue "D32" "3.519.29"
Set-TextValue "E32" "  -1.77%  "
Set-TextValue "D33" "0.153"
Set-TextValue "E33" "  -1.15%  "
Set-TextValue "D34" "24.11"
Set-TextValue "E34" "  -5.61%  "
Set-TextValue "D36" "1.32"
Set-TextValue "E36" "  -9.34%  "
Set-TextValue "D37" "7.50"
Set-TextValue "E37" "  -4.70%  "
Set-TextValue "D38" "173.63"
Set-TextValue "E38" "  -0.05%  "
Set-TextValue "D39" "5.21"
Set-TextValue "E39" "  -7.47%  "
Set-TextValue "E40" "  -8.97%  "
Set-TextValue "E41" "  -5.67%  "
Set-TextValue "D42" "4.98"
Set-TextValue "E42" "  -4.86%  "
Set-TextValue "D43" "0.854"
Set-TextValue "E43" "  -4.70%  "
Set-TextValue "D44" "45.17"
Set-TextValue "E44" "  -2.06%  "
Set-TextValue "E45" "  -8.18%  "
Set-TextValue "D46" "0.999"
Set-TextValue "E46" "  -0.06%  "
Set-TextValue "E47" "  -8.55%  "
Set-TextValue "D48" "23.17"
Set-TextValue "E48" "  -2.50%  "
Set-TextValue "E49" "  -2.93%  "
Set-TextValue "E50" "  -8.39%  "
Set-TextValue "D51" "0.900"
Set-TextValue "E51" "  -5.19%  "
